# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File", "Latest Handback File" and
#    "Latest Handback DateTime" columns populated for both data rows
#  - Column widths widen to fit the new, longer values
#  - New hyperlinks are added to the newly-populated "Latest Target File" cells

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$targetFile = "32148b73-c6c7-4f2c-b423-3e776891f11f.md"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8ec63a5d0dfce38fd0cf9f123c79a2b7bfb966a/e2e/32148b73-c6c7-4f2c-b423-3e776891f11f.md"

# --- Overview sheet ---------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusNew
$ov.Range("F2").Value = $statusNew
$ov.Range("E3").Value = $statusNew
$ov.Range("F3").Value = $statusNew
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet --------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

$zh.Range("I2").Value = $targetFile
$zh.Range("J2").Value = "32148b73-c6c7-4f2c-b423-3e776891f11f.b32bc2ec4764a6f6d6612e33d75134f195191199.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-24 07:04:17"

$zh.Range("I3").Value = $targetFile
$zh.Range("J3").Value = "32148b73-c6c7-4f2c-b423-3e776891f11f.b32bc2ec4764a6f6d6612e33d75134f195191199.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-24 07:04:17"

$zh.Hyperlinks.Add($zh.Range("I2"), $targetUrl, "", "", $targetFile)
$zh.Hyperlinks.Add($zh.Range("I3"), $targetUrl, "", "", $targetFile)

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet --------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

$de.Range("I2").Value = $targetFile
$de.Range("J2").Value = "32148b73-c6c7-4f2c-b423-3e776891f11f.b32bc2ec4764a6f6d6612e33d75134f195191199.de-de.xlf"
$de.Range("K2").Value = "2016-08-24 07:04:25"

$de.Range("I3").Value = $targetFile
$de.Range("J3").Value = "32148b73-c6c7-4f2c-b423-3e776891f11f.b32bc2ec4764a6f6d6612e33d75134f195191199.de-de.xlf"
$de.Range("K3").Value = "2016-08-24 07:04:25"

$de.Hyperlinks.Add($de.Range("I2"), $targetUrl, "", "", $targetFile)
$de.Hyperlinks.Add($de.Range("I3"), $targetUrl, "", "", $targetFile)

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
